# Remove the trailing page-break paragraph and the "[ANCIENTS]" suit-key
# paragraph that follows the "battling normally" effect text, then append a
# single trailing space run to the paragraph that now ends the document
# (mirrors the commit "Removed ancients from effect list").

$d = $word.ActiveDocument

# Locate the paragraph that ends with "...battling normally." (the anchor
# paragraph which should keep its trailing space run) and the paragraph
# that starts the "[ANCIENTS]" block (the end of the range to be removed).
$anchorIndex = -1
$ancientsIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "battling normally") {
        $anchorIndex = $i
    }
    if ($t -match "\[ANCIENTS\]") {
        $ancientsIndex = $i
    }
}

if ($anchorIndex -gt 0 -and $ancientsIndex -gt $anchorIndex) {
    $anchorPara = $d.Paragraphs.Item($anchorIndex)
    $lastRemoved = $d.Paragraphs.Item($ancientsIndex)

    # Delete everything from just after the anchor paragraph's own mark
    # through the end of the "[ANCIENTS]" paragraph -- this removes both
    # the intervening page-break paragraph and the whole ANCIENTS block,
    # paragraph marks included.
    $killRange = $d.Range($anchorPara.Range.End, $lastRemoved.Range.End)
    $killRange.Delete()

    # Re-fetch the anchor paragraph's range and append the single space
    # run the diff leaves behind.
    $anchorPara = $d.Paragraphs.Item($anchorIndex)
    $anchorPara.Range.InsertAfter(" ")
}
